$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.056318521499634
$ws.Range("B1").Value = 3.670591115951538
$ws.Range("C1").Value = 3.198753595352173
$ws.Range("D1").Value = 2.044086694717407
$ws.Range("E1").Value = 1.167993903160095
